$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1363.5454
$ws.Range("I15").Value = 1363.5454
$ws.Range("K15").Value = 4090.6362
$ws.Range("M15").Value = -3921.6362

$ws.Range("H62").Value = 983761.0600000001
$ws.Range("I62").Value = 1289459
$ws.Range("J62").Value = 168566.67
$ws.Range("K62").Value = 1289459
$ws.Range("L62").Value = 168566.67
$ws.Range("M62").Value = -1288835
$ws.Range("N62").Value = -169814.67

$ws.Range("H65").Value = 983761.0600000001
$ws.Range("I65").Value = 1289459
$ws.Range("J65").Value = 168566.67
$ws.Range("K65").Value = 6447295
$ws.Range("L65").Value = 842833.3500000001
$ws.Range("M65").Value = -6444175
$ws.Range("N65").Value = -849073.3500000001

$ws.Range("H128").Value = 49999
$ws.Range("J128").Value = 49999
$ws.Range("L128").Value = 49999
$ws.Range("N128").Value = -59959

$ws.Range("H132").Value = 3313.3333
$ws.Range("I132").Value = 3167
$ws.Range("K132").Value = 9501
$ws.Range("M132").Value = -6971

$ws.Range("H137").Value = 1279.2
$ws.Range("I137").Value = 1279.2
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 3837.6
$ws.Range("L137").Value = 0
$ws.Range("M137").Value = -1287.6
$ws.Range("N137").Value = $null

$ws.Range("H138").Value = 2133.7632
$ws.Range("I138").Value = 1204.24
$ws.Range("J138").Value = 3921.3076
$ws.Range("K138").Value = 3612.72
$ws.Range("L138").Value = 11763.9228
$ws.Range("M138").Value = 1527.28
$ws.Range("N138").Value = -22043.9228

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4783.5557
$ws.Range("I32").Value = 3800.2778
$ws.Range("J32").Value = 6750.1113
$ws.Range("K32").Value = 3800.2778
$ws.Range("L32").Value = 6750.1113
$ws.Range("M32").Value = -3513.2778
$ws.Range("N32").Value = -7324.1113

$ws.Range("H61").Value = 20835408
$ws.Range("I61").Value = 33334300
$ws.Range("K61").Value = 33334300
$ws.Range("M61").Value = -33334088

$ws.Range("H88").Value = 12822829
$ws.Range("I88").Value = 33334534
$ws.Range("J88").Value = 3014
$ws.Range("K88").Value = 33334534
$ws.Range("L88").Value = 3014
$ws.Range("M88").Value = -33334128
$ws.Range("N88").Value = -3826

$ws.Range("H91").Value = 12822829
$ws.Range("I91").Value = 33334534
$ws.Range("J91").Value = 3014
$ws.Range("K91").Value = 33334534
$ws.Range("L91").Value = 3014
$ws.Range("M91").Value = -33333130
$ws.Range("N91").Value = -5822

$ws.Range("H122").Value = 20836492
$ws.Range("I122").Value = 37038988
$ws.Range("K122").Value = 111116964
$ws.Range("M122").Value = -111114514

$ws.Range("H132").Value = 23811170
$ws.Range("I132").Value = 25642590
$ws.Range("K132").Value = 76927770
$ws.Range("M132").Value = -76925240

$ws.Range("H136").Value = 20835408
$ws.Range("I136").Value = 33334300
$ws.Range("K136").Value = 100002900
$ws.Range("M136").Value = -100000350

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H96").Value = 13038.077
$ws.Range("I96").Value = 13038.077
$ws.Range("J96").Value = 0
$ws.Range("K96").Value = 13038.077
$ws.Range("L96").Value = 0
$ws.Range("M96").Value = -10292.077
$ws.Range("N96").Value = $null

$ws.Range("H105").Value = 2111.3333
$ws.Range("I105").Value = 2015.2222
$ws.Range("J105").Value = 2399.6667
$ws.Range("K105").Value = 2015.2222
$ws.Range("L105").Value = 2399.6667
$ws.Range("M105").Value = -268.2221999999999
$ws.Range("N105").Value = -5893.6667

$ws.Range("H134").Value = 2087.6667
$ws.Range("I134").Value = 632
$ws.Range("J134").Value = 4999
$ws.Range("K134").Value = 1896
$ws.Range("L134").Value = 14997
$ws.Range("M134").Value = 639
$ws.Range("N134").Value = -20067

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H20").Value = 55999
$ws.Range("J20").Value = 55999
$ws.Range("L20").Value = 55999
$ws.Range("N20").Value = -56471

$ws.Range("H30").Value = 55999
$ws.Range("J30").Value = 55999
$ws.Range("L30").Value = 55999
$ws.Range("N30").Value = -56181

$ws.Range("H31").Value = 3287.0667
$ws.Range("I31").Value = 1703.8889
$ws.Range("K31").Value = 1703.8889
$ws.Range("M31").Value = -1408.8889

$ws.Range("H34").Value = 3287.0667
$ws.Range("I34").Value = 1703.8889
$ws.Range("K34").Value = 1703.8889
$ws.Range("M34").Value = -1501.8889

$ws.Range("H122").Value = 2300.6365
$ws.Range("I122").Value = 1162
$ws.Range("J122").Value = 3249.5
$ws.Range("K122").Value = 3486
$ws.Range("L122").Value = 9748.5
$ws.Range("M122").Value = -1036
$ws.Range("N122").Value = -14648.5

$ws.Range("H128").Value = 55999
$ws.Range("J128").Value = 55999
$ws.Range("L128").Value = 55999
$ws.Range("N128").Value = -65959

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 1722
$ws.Range("I3").Value = 1359.8
$ws.Range("J3").Value = 3533
$ws.Range("K3").Value = 4079.4
$ws.Range("L3").Value = 10599
$ws.Range("M3").Value = -3967.4
$ws.Range("N3").Value = -10823

$ws.Range("H5").Value = 1332.4286
$ws.Range("J5").Value = 1192.6666
$ws.Range("L5").Value = 3577.9998
$ws.Range("N5").Value = -3801.9998

$ws.Range("H18").Value = 1501.6666
$ws.Range("I18").Value = 1501.6666
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 4504.9998
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = -4335.9998
$ws.Range("N18").Value = $null

$ws.Range("H22").Value = 1133.3334
$ws.Range("I22").Value = 700
$ws.Range("K22").Value = 2100
$ws.Range("M22").Value = -1931

$ws.Range("H27").Value = 1133.3334
$ws.Range("I27").Value = 700
$ws.Range("K27").Value = 2100
$ws.Range("M27").Value = -1998

$ws.Range("H74").Value = 7343.3335
$ws.Range("I74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("M74").Value = $null

$ws.Range("H77").Value = 7343.3335
$ws.Range("I77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("M77").Value = $null

$ws.Range("H80").Value = 5913
$ws.Range("J80").Value = 6428.5713
$ws.Range("L80").Value = 19285.7139
$ws.Range("N80").Value = -21157.7139

$ws.Range("H83").Value = 5913
$ws.Range("J83").Value = 6428.5713
$ws.Range("L83").Value = 57857.14169999999
$ws.Range("N83").Value = -67217.14169999999

$ws.Range("H126").Value = 1000
$ws.Range("I126").Value = 1000
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 3000
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = 1940
$ws.Range("N126").Value = $null

$ws.Range("H135").Value = 1332.4286
$ws.Range("J135").Value = 1192.6666
$ws.Range("L135").Value = 10733.9994
$ws.Range("N135").Value = -15803.9994

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2304.25
$ws.Range("I122").Value = 2304.25
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 6912.75
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -4462.75
$ws.Range("N122").Value = $null

$ws.Range("H132").Value = 3998.7932
$ws.Range("I132").Value = 4445.1113
$ws.Range("J132").Value = 3268.4546
$ws.Range("K132").Value = 13335.3339
$ws.Range("L132").Value = 9805.363799999999
$ws.Range("M132").Value = -10805.3339
$ws.Range("N132").Value = -14865.3638

$ws.Range("H135").Value = 44995
$ws.Range("J135").Value = 44995
$ws.Range("L135").Value = 44995
$ws.Range("N135").Value = -55135

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 18520968
$ws.Range("I7").Value = 27779588
$ws.Range("J7").Value = 3727.2222
$ws.Range("K7").Value = 27779588
$ws.Range("L7").Value = 3727.2222
$ws.Range("M7").Value = -27779476
$ws.Range("N7").Value = -3951.2222

$ws.Range("H68").Value = 4374.75
$ws.Range("I68").Value = 5766.3335
$ws.Range("J68").Value = 3539.8
$ws.Range("K68").Value = 5766.3335
$ws.Range("L68").Value = 3539.8
$ws.Range("M68").Value = -5017.3335
$ws.Range("N68").Value = -5037.8

$ws.Range("H71").Value = 4374.75
$ws.Range("I71").Value = 5766.3335
$ws.Range("J71").Value = 3539.8
$ws.Range("K71").Value = 28831.6675
$ws.Range("L71").Value = 17699
$ws.Range("M71").Value = -25087.6675
$ws.Range("N71").Value = -25187

$ws.Range("H82").Value = 1170.25
$ws.Range("I82").Value = 1397
$ws.Range("J82").Value = 490
$ws.Range("K82").Value = 1397
$ws.Range("L82").Value = 490
$ws.Range("M82").Value = -1036
$ws.Range("N82").Value = -1212

$ws.Range("H85").Value = 1170.25
$ws.Range("I85").Value = 1397
$ws.Range("J85").Value = 490
$ws.Range("K85").Value = 1397
$ws.Range("L85").Value = 490
$ws.Range("M85").Value = -149
$ws.Range("N85").Value = -2986

$ws.Range("H122").Value = 4895.6665
$ws.Range("I122").Value = 4928.5
$ws.Range("J122").Value = 4849.7
$ws.Range("K122").Value = 14785.5
$ws.Range("L122").Value = 14549.1
$ws.Range("M122").Value = -12335.5
$ws.Range("N122").Value = -19449.1

$ws.Range("H126").Value = 18520968
$ws.Range("I126").Value = 27779588
$ws.Range("J126").Value = 3727.2222
$ws.Range("K126").Value = 83338764
$ws.Range("L126").Value = 11181.6666
$ws.Range("M126").Value = -83336294
$ws.Range("N126").Value = -16121.6666

$ws.Range("H132").Value = 5434
$ws.Range("I132").Value = 2592.6
$ws.Range("K132").Value = 7777.799999999999
$ws.Range("M132").Value = -5247.799999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H49").Value = 50000
$ws.Range("I49").Value = 50000
$ws.Range("K49").Value = 50000
$ws.Range("M49").Value = -49770

$ws.Range("H122").Value = 1856.7778
$ws.Range("I122").Value = 1837.3572
$ws.Range("J122").Value = 1924.75
$ws.Range("K122").Value = 5512.071599999999
$ws.Range("L122").Value = 5774.25
$ws.Range("M122").Value = -3062.071599999999
$ws.Range("N122").Value = -10674.25

$ws.Range("H132").Value = 3305.3914
$ws.Range("I132").Value = 2745
$ws.Range("K132").Value = 8235
$ws.Range("M132").Value = -5705
